$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): update F column "想去人数" counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 874
$ws1.Range("F8").Value = 4598
$ws1.Range("F16").Value = 2595
$ws1.Range("F21").Value = 2293
$ws1.Range("F25").Value = 161
$ws1.Range("F26").Value = 106

# Sheet "全部类型" (sheet4.xml): same rows, offset by one because of an extra entry
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 874
$ws4.Range("F9").Value = 4598
$ws4.Range("F17").Value = 2595
$ws4.Range("F22").Value = 2293
$ws4.Range("F26").Value = 161
$ws4.Range("F27").Value = 106
